# Adapt RestDataProvider to multiple scenarii counters
#
# Fills in the "Result" column (F) of the NoraUi-blog worksheet with the
# outcome of each scenario run: a success, a failure (anonymous users are
# not allowed to post) or an ignored step (skipped because of an earlier
# failure). Each outcome gets its own font color so the column reads like
# a simple status report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Article 3 -> succeeded
$ws.Range("F4").ClearFormats()
$ws.Range("F4").Value = "Succès"
$ws.Range("F4").Font.ColorIndex = 17

# Row 2: Article 1 -> rejected because posted by "anonymous"
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "Échec : anonymous is prohibited in demo blog!!"
$ws.Range("F2").Font.ColorIndex = 10

# Row 3: Article 2 -> skipped, a previous step already failed
$ws.Range("F3").ClearFormats()
$ws.Range("F3").Value = "Élément ignoré suite à une erreur précédente."
$ws.Range("F3").Font.ColorIndex = 53

# Row 6: Article 4 -> succeeded (also drops the row's old highlight fill)
$ws.Range("F6").ClearFormats()
$ws.Range("F6").Value = "Succès"
$ws.Range("F6").Font.ColorIndex = 17

# Row 7: Article 5 -> succeeded
$ws.Range("F7").ClearFormats()
$ws.Range("F7").Value = "Succès"
$ws.Range("F7").Font.ColorIndex = 17
